$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text '2024-02-19 Monday' '2024-02-20 Tuesday'
Replace-Text '69×45=' '91×97='
Replace-Text '20×34=' '86×43='
Replace-Text '18×41=' '22×47='
Replace-Text '25×28=' '65×82='
Replace-Text '37×82=' '82×15='
Replace-Text '76×17=' '21×84='
Replace-Text '42×11=' '80×23='
Replace-Text '20×80=' '56×44='
Replace-Text '13×11=' '92×98='
Replace-Text '82×84=' '25×36='
Replace-Text '31×94=' '98×41='
Replace-Text '77×70=' '59×80='
Replace-Text '59×12=' '19×35='
Replace-Text '87×98=' '52×94='
Replace-Text '26×35=' '88×65='
Replace-Text '35×63=' '67×91='
Replace-Text '69×35=' '15×23='
Replace-Text '47×14=' '43×98='
Replace-Text '14×40=' '94×45='
Replace-Text '75×33=' '94×94='
Replace-Text '57×66=' '35×65='
Replace-Text '41×49=' '36×30='
Replace-Text '12×97=' '62×17='
Replace-Text '44×19=' '83×22='
Replace-Text '67×74=' '88×78='
